# ES001D.xlsx "View Dfn" sheet update
# 1. Django bump note changed rev/timestamp/class path fields.
# 2. Minor table bookkeeping cells added (blank filler cells in col A/B).
# 3. The unused "schClaimer" search-field row (old row 36) was removed,
#    shifting the rest of the fieldTable block (and everything below it)
#    up by one row.
# 4. Selection cursor moved to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("View Dfn")

# --- simple value edits -----------------------------------------------
$ws.Range("C10").Value = "es.views.es001.ES001D"

# C14 is formatted as Text ("@") even though it stores a number, so a
# straight .Value assignment would get coerced to a string. Flip to
# General just long enough to write a real numeric value, then restore
# the original (text) number format to match the original cell style.
$c14 = $ws.Range("C14")
$origFormat = $c14.NumberFormat
$c14.NumberFormat = "General"
$c14.Value2 = 20250821173900
$c14.NumberFormat = $origFormat

# --- newly-added blank filler cells (match surrounding column style) --
# Column A "blank" cells use the same style as A19 (etc.); column B
# "blank" cells use the same style as B27 (etc.). Copy+PasteSpecial
# (formats only) reproduces the exact style index without touching value.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)

$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("B29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- remove the obsolete "schClaimer" row; everything below shifts up -
$ws.Rows.Item(36).Delete()

# --- move the active selection ----------------------------------------
$ws.Range("C3").Select() | Out-Null
